# The pickle-load bug meant the clustering object wasn't being received,
# so the centroid assignment / distance columns (B:H, I, J) for the data
# points were computed against stale/incorrect state. After the fix, the
# per-row cluster one-hot flags and the Edad_Escalada / Edad values are
# recomputed (rows essentially re-shuffle between each other).
#
# Apply the corrected values cell-by-cell, row by row (rows 2-11).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (A2 = 0)
$ws.Range("B2").Value = -0.0
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 0
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = -0.0
$ws.Range("I2").Value = 0.5604
$ws.Range("J2").Value = 44.2716

# Row 3 (A3 = 1)
$ws.Range("E3").Value = -0.0

# Row 4 (A4 = 2)
$ws.Range("B4").Value = 0
$ws.Range("D4").Value = -0.0
$ws.Range("E4").Value = 1
$ws.Range("G4").Value = -0.0
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0.5713
$ws.Range("J4").Value = 45.1327

# Row 5 (A5 = 3)
$ws.Range("H5").Value = 0

# Row 6 (A6 = 4)
$ws.Range("B6").Value = 1
$ws.Range("C6").Value = -0.0
$ws.Range("F6").Value = 0
$ws.Range("I6").Value = 0.5661
$ws.Range("J6").Value = 44.72190000000001

# Row 7 (A7 = 5)
$ws.Range("B7").Value = 0
$ws.Range("C7").Value = 1
$ws.Range("I7").Value = 0.5339
$ws.Range("J7").Value = 42.1781

# Row 8 (A8 = 6)
$ws.Range("B8").Value = 0
$ws.Range("C8").Value = 1
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = 1
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0.5513
$ws.Range("J8").Value = 43.5527

# Row 9 (A9 = 7)
$ws.Range("B9").Value = 1
$ws.Range("C9").Value = 0
$ws.Range("I9").Value = 0.5456
$ws.Range("J9").Value = 43.1024

# Row 10 (A10 = 8)
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = -0.0
$ws.Range("I10").Value = 0.5737
$ws.Range("J10").Value = 45.3223

# Row 11 (A11 = 9)
$ws.Range("G11").Value = 0
